$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999706412279,
    0.99776509643768041,
    0.99527213430868566,
    1.0007870180457681,
    0.99325365079342465,
    0.97552217496641414,
    0.97070409955833004,
    0.96528693306256841,
    0.95942533995366253,
    0.95449174875873943,
    0.95390021194060837,
    0.95303065587640035,
    0.95033960705251452,
    0.95030795386674805,
    0.95136319711425954,
    0.95305324530841484,
    0.94934558777402533,
    0.94823667362179465,
    0.99417484375844256,
    0.96466130303983588,
    0.96326276768106267,
    0.96199825039250009,
    0.98630730615018636,
    0.96933781588546331,
    0.9628808617129796,
    0.9742483594511806,
    0.97272422564491956,
    0.96712288070173469,
    0.96367456493388048,
    0.96285423015369187,
    0.95520200370661279,
    0.95352270838591213,
    0.95300269184011088
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

